# Refresh the crypto "Price" / "Volume(1h)" columns (and, for rows 10-18,
# the Coin/Link columns too, since the ranking reshuffled) to match the
# latest scrape. Price values are stored as plain text in this sheet (not
# numbers), so numeric-looking strings are written with a leading "'" to
# keep Excel from silently coercing them into Number cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.41"
$ws.Range("D3").Value = "'24.00"
$ws.Range("D4").Value = "'5.350"
$ws.Range("D5").Value = "'0.05853"
$ws.Range("D6").Value = "'6.487"
$ws.Range("D7").Value = "'3.364"
$ws.Range("D8").Value = "'0.8121"
$ws.Range("D9").Value = "'0.9257"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01078"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1412"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07377"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03061"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03060"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09338"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.861"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001560"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04694"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").Value = "'0.005983"
$ws.Range("D20").Value = "'0.001246"
$ws.Range("D22").Value = "'0.00008814"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("D23").Value = "'3.602"
$ws.Range("D26").Value = "'0.1333"
$ws.Range("D28").Value = "'0.0002657"
$ws.Range("D40").Value = "'0.03845"
$ws.Range("D41").Value = "'0.006411"
$ws.Range("D42").Value = "'0.1063"
$ws.Range("D43").Value = "'0.002950"
$ws.Range("D44").Value = "'0.008260"
$ws.Range("D45").Value = "'0.00005271"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D47").Value = "'0.6531"
$ws.Range("D48").Value = "'0.001729"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
$ws.Range("D49").Value = "'0.00002104"
$ws.Range("D50").Value = "'0.0002004"
